# Generate Report for Handback
#
# Row 7 in both the "zh-cn" and "de-de" sheets represents the
# 0faceb12-c00f-4929-922b-d73e31f6d63b.md handback file. A handback was
# received for this file, but its base version is out of date, so we
# record the (now available) target/handback file info, the handback
# timestamp, and a "stale version" error message.

$wb = $excel.ActiveWorkbook

$sourceFile   = "0faceb12-c00f-4929-922b-d73e31f6d63b.md"
$currentUrl   = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aad5f3ce6cf12d532bb3fa6e12903fc99a9d393a/e2e/0faceb12-c00f-4929-922b-d73e31f6d63b.md"
$latestUrl    = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d8d209cbb5b1d7fe099e706b9440814712810fba/e2e/0faceb12-c00f-4929-922b-d73e31f6d63b.md"
$errorDetail  = "The version of handback file is not the latest, current: $currentUrl, latest: $latestUrl."

function Set-HandbackRow {
    param(
        $Workbook,
        [string]$SheetName,
        [string]$HandbackFile,
        [string]$HandbackDateTime,
        [string]$SourceFile,
        [string]$CurrentUrl,
        [string]$ErrorDetail
    )

    $ws = $Workbook.Worksheets.Item($SheetName)

    # I7 - Latest Target File: the handed-back source file is now available.
    $ws.Range("I7").Value = $SourceFile
    $ws.Hyperlinks.Add($ws.Range("I7"), $CurrentUrl, "", "", $SourceFile)
    $ws.Range("I7").Font.Underline = 2
    $ws.Range("I7").Font.Color = 15570276

    # J7 - Latest Handback File
    $ws.Range("J7").Value = $HandbackFile

    # K7 - Latest Handback DateTime
    $ws.Range("K7").Value = $HandbackDateTime

    # P7 - Error Detail
    $ws.Range("P7").Value = $ErrorDetail
}

Set-HandbackRow $wb "zh-cn" "0faceb12-c00f-4929-922b-d73e31f6d63b.4b3ee365de5f190a8fb1be3df3a4062cd9805c53.zh-cn.xlf" "2016-09-02 03:03:45" $sourceFile $currentUrl $errorDetail
Set-HandbackRow $wb "de-de" "0faceb12-c00f-4929-922b-d73e31f6d63b.4b3ee365de5f190a8fb1be3df3a4062cd9805c53.de-de.xlf" "2016-09-02 03:03:52" $sourceFile $currentUrl $errorDetail
